$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.28991759397096928,
    -0.0099999994008044268,
    -0.0089999993921114907,
    -0.011999999847677856,
    -0.0059999994013333335,
    0.071913365167972643,
    -0.019999999282211078,
    -0.090532356016239213,
    -0.005999999379238119,
    -0.0059999993824746412,
    0.080358353884715683,
    -0.0059999993796391315,
    -0.0059999993698625076,
    -0.011999999323326627,
    -0.0059999993649917371,
    -0.0059999993630666104,
    -0.0059999993604478163,
    -0.0089999993379983323,
    -0.0089999994074987377,
    -0.0089999994021727758,
    -0.0089999994013965079,
    -0.0089999994008653772,
    -0.0089999993867060368,
    -0.041999999137670052,
    -0.041999999132922738,
    -0.0059999993848052213,
    -0.0059999993817090314,
    -0.005999999367853448,
    -0.011999999315358778,
    -0.019999999252570344,
    -0.01499999928263307,
    -0.02099999923789575,
    -0.0059999993469173063
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
